$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Copy the formatting of the last existing data row (21) down onto the new
# row (24) before writing any values into it, so the new cells pick up the
# same styles (s="5" for B/C/D, s="7" for E:H) as the rest of the table.
$ws.Range("B21:H21").Copy()
$ws.Range("B24").PasteSpecial(-4122)  # xlPasteFormats

# New test case fm19: "Residential policy with blanket policy terms.
# Previous level input loss back-allocation", allocrule 2.
$ws.Range("B24").Value = "fm19"

# Existing fm17/fm18 rows: allocrule changed from -1 to 0 (still stored as
# text, keeping the leading apostrophe so Excel treats "0,1"/"0,2" as text
# rather than trying to interpret them as numbers).
$ws.Range("D22").Value = "'0,1"
$ws.Range("D23").Value = "'0,2"

$ws.Range("C24").Value = "Residential policy with blanket policy terms. Previous level input loss back-allocation"
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 1

# Match the saved cursor position after entering the new row of data.
$null = $ws.Range("I24").Select()
